# Actualiza base de datos EC:
# El periodo en mora reportado cambia de 2507 a 2508 para los
# trabajadores listados en la hoja de Estado de Cuenta.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
